# AFDP-1050 - Implement Document level security - initial implementation of folder security
#
# Adds a new "Folder - default public access" rule row to the rules table on
# Sheet1 (row 39), tags the sheets with a white tab color, and updates the
# active selection on Sheet1.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# New rule row describing default public access for folders.
$ws1.Range("B39").Value = "Folder – default public access"
$ws1.Range("C39").Value = "FOLDER"
$ws1.Range("G39").Value = "grant read to *"

# Row 39 grows a bit taller to fit the wrapped text, matching the other rows.
$ws1.Rows.Item(39).RowHeight = 23.85

# Tag every sheet with a (white) tab color.
$ws1.Tab.Color = 16777215
$ws2.Tab.Color = 16777215
$ws3.Tab.Color = 16777215

# Update the remembered selection on the rules sheet to the new last row.
$null = $ws1.Activate()
$null = $ws1.Range("G40").Select()
